$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers ("from" row) - add pv1 / bat1 alongside the existing net1
$ws.Range("B1").Value2 = "P_from_net1"
$ws.Range("C1").Value2 = "P_from_pv1"
$ws.Range("D1").Value2 = "P_from_bat1"

# Column A row labels - demand1/demand2/to_net1/to_bat1 (replaces the old
# charging_station1 rows with a second demand row + a battery row)
$ws.Range("A2").Value2 = "param_P_to_demand1"
$ws.Range("A3").Value2 = "param_P_to_demand2"
$ws.Range("A4").Value2 = "P_to_net1"
$ws.Range("A5").Value2 = "P_to_bat1"

# Column B data (net1)
$ws.Range("B2").Value2 = "P_net1_demand1"
$ws.Range("B3").Value2 = "P_net1_demand2"
$ws.Range("B4").Value2 = 0
$ws.Range("B5").Value2 = "P_net1_bat1"

# Column C data (pv1) - new column
$ws.Range("C2").Value2 = "P_pv1_demand1"
$ws.Range("C3").Value2 = "P_pv1_demand2"
$ws.Range("C4").Value2 = "P_pv1_net1"
$ws.Range("C5").Value2 = "P_pv1_bat1"

# Column D data (bat1) - new column
$ws.Range("D2").Value2 = "P_bat1_demand1"
$ws.Range("D3").Value2 = "P_bat1_demand2"
$ws.Range("D4").Value2 = "P_bat1_net1"
$ws.Range("D5").Value2 = 0

# Make sure every "label" cell (bold, centered, top-aligned, bordered) uses
# the same style as before - copy the format from the original B1 label
# cell onto the whole header row + row-label column (one area at a time -
# PasteSpecial onto a multi-area range only fills the first area).
$ws.Range("B1").Copy()
foreach ($addr in @("B1","C1","D1","A2","A3","A4","A5")) {
    $ws.Range($addr).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
